$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.020574846168834
$ws.Cells.Item(2, 4).Value = 1.026259465027548
$ws.Cells.Item(2, 5).Value = 1.02155833181374
$ws.Cells.Item(2, 6).Value = 1.031506699814445
$ws.Cells.Item(2, 9).Value = 1.029430903032845
$ws.Cells.Item(2, 10).Value = 1.025771184298331
$ws.Cells.Item(2, 11).Value = 1.029082633784605
$ws.Cells.Item(2, 12).Value = 1.024395307221931
$ws.Cells.Item(2, 13).Value = 1.034314619661036
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.021560697842817
$ws.Cells.Item(3, 4).Value = 1.026984683481958
$ws.Cells.Item(3, 5).Value = 1.022395848360031
$ws.Cells.Item(3, 6).Value = 1.03270218394147
$ws.Cells.Item(3, 9).Value = 1.029600306709796
$ws.Cells.Item(3, 10).Value = 1.026394092810483
$ws.Cells.Item(3, 11).Value = 1.029615551789267
$ws.Cells.Item(3, 12).Value = 1.025039217348184
$ws.Cells.Item(3, 13).Value = 1.035317645065252
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.022198891291062
$ws.Cells.Item(4, 4).Value = 1.027453956667276
$ws.Cells.Item(4, 5).Value = 1.022938399991487
$ws.Cells.Item(4, 6).Value = 1.033476155489049
$ws.Cells.Item(4, 9).Value = 1.029708574408544
$ws.Cells.Item(4, 10).Value = 1.026796870787447
$ws.Cells.Item(4, 11).Value = 1.029959727298783
$ws.Cells.Item(4, 12).Value = 1.02545585119069
$ws.Cells.Item(4, 13).Value = 1.035966519531226
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.022467254880221
$ws.Cells.Item(5, 4).Value = 1.027651240124347
$ws.Cells.Item(5, 5).Value = 1.02316663682961
$ws.Cells.Item(5, 6).Value = 1.033801632493445
$ws.Cells.Item(5, 9).Value = 1.029753767138342
$ws.Cells.Item(5, 10).Value = 1.026966129773703
$ws.Cells.Item(5, 11).Value = 1.030104260478778
$ws.Cells.Item(5, 12).Value = 1.025630998866819
$ws.Cells.Item(5, 13).Value = 1.036239270636496
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.022512318230924
$ws.Cells.Item(6, 4).Value = 1.027684364915302
$ws.Cells.Item(6, 5).Value = 1.023204967460064
$ws.Cells.Item(6, 6).Value = 1.033856287353314
$ws.Cells.Item(6, 9).Value = 1.029761336245072
$ws.Cells.Item(6, 10).Value = 1.026994545047993
$ws.Cells.Item(6, 11).Value = 1.030128518939313
$ws.Cells.Item(6, 12).Value = 1.025660406609689
$ws.Cells.Item(6, 13).Value = 1.036285064658478
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022202476915362
$ws.Cells.Item(7, 4).Value = 1.027456592775766
$ws.Cells.Item(7, 5).Value = 1.022941449122333
$ws.Cells.Item(7, 6).Value = 1.033480504137763
$ws.Cells.Item(7, 9).Value = 1.02970917954565
$ws.Cells.Item(7, 10).Value = 1.026799132705604
$ws.Cells.Item(7, 11).Value = 1.029961659180674
$ws.Cells.Item(7, 12).Value = 1.025458191543697
$ws.Cells.Item(7, 13).Value = 1.03597016418495
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.02090796087583
$ws.Cells.Item(8, 4).Value = 1.02650455367644
$ws.Cells.Item(8, 5).Value = 1.021841245048126
$ws.Cells.Item(8, 6).Value = 1.031910633920977
$ws.Cells.Item(8, 9).Value = 1.02948843266502
$ws.Cells.Item(8, 10).Value = 1.025981757939122
$ws.Cells.Item(8, 11).Value = 1.029262871781822
$ws.Cells.Item(8, 12).Value = 1.024612923235103
$ws.Cells.Item(8, 13).Value = 1.034653627979123
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018629028235771
$ws.Cells.Item(9, 4).Value = 1.024827046526448
$ws.Cells.Item(9, 5).Value = 1.019907350405284
$ws.Cells.Item(9, 6).Value = 1.029147469659692
$ws.Cells.Item(9, 9).Value = 1.029089140896511
$ws.Cells.Item(9, 10).Value = 1.024539275869282
$ws.Cells.Item(9, 11).Value = 1.028026507940745
$ws.Cells.Item(9, 12).Value = 1.023123333902997
$ws.Cells.Item(9, 13).Value = 1.032332564516079
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.017111215426555
$ws.Cells.Item(10, 4).Value = 1.023708837850877
$ws.Cells.Item(10, 5).Value = 1.018621368710492
$ws.Cells.Item(10, 6).Value = 1.027307447864711
$ws.Cells.Item(10, 9).Value = 1.028816035772146
$ws.Cells.Item(10, 10).Value = 1.0235761964207
$ws.Cells.Item(10, 11).Value = 1.02719893386481
$ws.Cells.Item(10, 12).Value = 1.022130230277165
$ws.Cells.Item(10, 13).Value = 1.030784394130465
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.016454338434648
$ws.Cells.Item(11, 4).Value = 1.023224683492675
$ws.Cells.Item(11, 5).Value = 1.018065312921483
$ws.Cells.Item(11, 6).Value = 1.02651118600082
$ws.Cells.Item(11, 9).Value = 1.028696143035387
$ws.Cells.Item(11, 10).Value = 1.023158840244996
$ws.Cells.Item(11, 11).Value = 1.02683980249817
$ws.Cells.Item(11, 12).Value = 1.021700202067716
$ws.Cells.Item(11, 13).Value = 1.030113826331614
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016210397172966
$ws.Cells.Item(12, 4).Value = 1.023044853510665
$ws.Cells.Item(12, 5).Value = 1.017858887415823
$ws.Cells.Item(12, 6).Value = 1.026215490224437
$ws.Cells.Item(12, 9).Value = 1.028651363952704
$ws.Cells.Item(12, 10).Value = 1.023003765438975
$ws.Cells.Item(12, 11).Value = 1.02670628753073
$ws.Cells.Item(12, 12).Value = 1.021540469790147
$ws.Cells.Item(12, 13).Value = 1.029864717084059
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.01626272104202
$ws.Cells.Item(13, 4).Value = 1.023083427343721
$ws.Cells.Item(13, 5).Value = 1.017903161018545
$ws.Cells.Item(13, 6).Value = 1.026278914756068
$ws.Cells.Item(13, 9).Value = 1.028660980326695
$ws.Cells.Item(13, 10).Value = 1.023037031773723
$ws.Cells.Item(13, 11).Value = 1.026734932256741
$ws.Cells.Item(13, 12).Value = 1.021574732927204
$ws.Cells.Item(13, 13).Value = 1.029918153251786
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.016434173110304
$ws.Cells.Item(14, 4).Value = 1.023209818552466
$ws.Cells.Item(14, 5).Value = 1.018048247305572
$ws.Cells.Item(14, 6).Value = 1.026486742244779
$ws.Cells.Item(14, 9).Value = 1.028692446590673
$ws.Cells.Item(14, 10).Value = 1.023146022733652
$ws.Cells.Item(14, 11).Value = 1.026828768504885
$ws.Cells.Item(14, 12).Value = 1.021686998558172
$ws.Cells.Item(14, 13).Value = 1.030093235508729
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.016539817295755
$ws.Cells.Item(15, 4).Value = 1.02328769323592
$ws.Cells.Item(15, 5).Value = 1.01813765546402
$ws.Cells.Item(15, 6).Value = 1.026614801050094
$ws.Cells.Item(15, 9).Value = 1.028711801460543
$ws.Cells.Item(15, 10).Value = 1.023213169029312
$ws.Cells.Item(15, 11).Value = 1.026886568549292
$ws.Cells.Item(15, 12).Value = 1.021756169047472
$ws.Cells.Item(15, 13).Value = 1.030201105408151
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.017154817488713
$ws.Cells.Item(16, 4).Value = 1.023740970444469
$ws.Cells.Item(16, 5).Value = 1.018658288859747
$ws.Cells.Item(16, 6).Value = 1.027360303186138
$ws.Cells.Item(16, 9).Value = 1.028823958203942
$ws.Cells.Item(16, 10).Value = 1.023603887922038
$ws.Cells.Item(16, 11).Value = 1.027222751710706
$ws.Cells.Item(16, 12).Value = 1.022158769713194
$ws.Cells.Item(16, 13).Value = 1.030828893255118
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.01754068357344
$ws.Cells.Item(17, 4).Value = 1.024025309896808
$ws.Cells.Item(17, 5).Value = 1.018985078548248
$ws.Cells.Item(17, 6).Value = 1.027828064891659
$ws.Cells.Item(17, 9).Value = 1.028893873242118
$ws.Cells.Item(17, 10).Value = 1.023848885776897
$ws.Cells.Item(17, 11).Value = 1.027433420391468
$ws.Cells.Item(17, 12).Value = 1.022411308754401
$ws.Cells.Item(17, 13).Value = 1.0312226341444
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.01776578610511
$ws.Cells.Item(18, 4).Value = 1.024191163835686
$ws.Cells.Item(18, 5).Value = 1.01917576491919
$ws.Cells.Item(18, 6).Value = 1.028100948715533
$ws.Cells.Item(18, 9).Value = 1.028934495592433
$ws.Cells.Item(18, 10).Value = 1.023991756393982
$ws.Cells.Item(18, 11).Value = 1.02755622396241
$ws.Cells.Item(18, 12).Value = 1.022558609718706
$ws.Cells.Item(18, 13).Value = 1.031452277298019
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.017842545908157
$ws.Cells.Item(19, 4).Value = 1.024247716317362
$ws.Cells.Item(19, 5).Value = 1.019240796822228
$ws.Cells.Item(19, 6).Value = 1.028194002929247
$ws.Cells.Item(19, 9).Value = 1.028948319971845
$ws.Cells.Item(19, 10).Value = 1.024040466059724
$ws.Cells.Item(19, 11).Value = 1.027598083914805
$ws.Cells.Item(19, 12).Value = 1.022608835391934
$ws.Cells.Item(19, 13).Value = 1.031530576429512
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.017499280341464
$ws.Cells.Item(20, 4).Value = 1.023994802591184
$ws.Cells.Item(20, 5).Value = 1.018950009295077
$ws.Cells.Item(20, 6).Value = 1.027777873700029
$ws.Cells.Item(20, 9).Value = 1.028886388351433
$ws.Cells.Item(20, 10).Value = 1.023822603182974
$ws.Cells.Item(20, 11).Value = 1.027410825483091
$ws.Cells.Item(20, 12).Value = 1.02238421378278
$ws.Cells.Item(20, 13).Value = 1.031180391467362
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016383683310296
$ws.Cells.Item(21, 4).Value = 1.023172599305058
$ws.Cells.Item(21, 5).Value = 1.018005519736376
$ws.Cells.Item(21, 6).Value = 1.026425540266078
$ws.Cells.Item(21, 9).Value = 1.028683187337295
$ws.Cells.Item(21, 10).Value = 1.023113928993338
$ws.Cells.Item(21, 11).Value = 1.026801139310195
$ws.Cells.Item(21, 12).Value = 1.021653939141899
$ws.Cells.Item(21, 13).Value = 1.030041678994028
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015682564679571
$ws.Cells.Item(22, 4).Value = 1.022655685531479
$ws.Cells.Item(22, 5).Value = 1.017412366912827
$ws.Cells.Item(22, 6).Value = 1.025575686912428
$ws.Cells.Item(22, 9).Value = 1.028554006067373
$ws.Cells.Item(22, 10).Value = 1.022668067510752
$ws.Cells.Item(22, 11).Value = 1.026417125361156
$ws.Cells.Item(22, 12).Value = 1.021194782964415
$ws.Cells.Item(22, 13).Value = 1.029325548503906
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.0160542123109
$ws.Cells.Item(23, 4).Value = 1.022929707399185
$ws.Cells.Item(23, 5).Value = 1.017726743236917
$ws.Cells.Item(23, 6).Value = 1.02602617141469
$ws.Cells.Item(23, 9).Value = 1.028622622080673
$ws.Cells.Item(23, 10).Value = 1.022904454462379
$ws.Cells.Item(23, 11).Value = 1.026620762699266
$ws.Cells.Item(23, 12).Value = 1.021438190503691
$ws.Cells.Item(23, 13).Value = 1.029705199744881
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.01751798857659
$ws.Cells.Item(24, 4).Value = 1.024008587518509
$ws.Cells.Item(24, 5).Value = 1.01896585534787
$ws.Cells.Item(24, 6).Value = 1.027800552795604
$ws.Cells.Item(24, 9).Value = 1.028889770939256
$ws.Cells.Item(24, 10).Value = 1.023834479256197
$ws.Cells.Item(24, 11).Value = 1.027421035384137
$ws.Cells.Item(24, 12).Value = 1.022396456836589
$ws.Cells.Item(24, 13).Value = 1.03119947917402
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.01921792808667
$ws.Cells.Item(25, 4).Value = 1.025260703701646
$ws.Cells.Item(25, 5).Value = 1.02040673384436
$ws.Cells.Item(25, 6).Value = 1.029861443796837
$ws.Cells.Item(25, 9).Value = 1.029193586869212
$ws.Cells.Item(25, 10).Value = 1.024912445355594
$ws.Cells.Item(25, 11).Value = 1.028346727587479
$ws.Cells.Item(25, 12).Value = 1.023508438825193
$ws.Cells.Item(25, 13).Value = 1.032932753981723
